{"js": "// Update the two-digit-divided-by-one-digit answer table: each table\n// cell holds a single \"XX\u00f7Y=Q, R\" style answer string that is being\n// replaced with a new equation/answer string (same cell, same\n// formatting - only the w:t text content changes).\nconst pairs = [\n  [\"74\u00f78=9, 2\", \"37\u00f76=6, 1\"],\n  [\"94\u00f78=11, 6\", \"75\u00f73=25, 0\"],\n  [\"65\u00f72=32, 1\", \"84\u00f72=42, 0\"],\n  [\"94\u00f76=15, 4\", \"36\u00f72=18, 0\"],\n  [\"25\u00f73=8, 1\", \"86\u00f74=21, 2\"],\n  [\"54\u00f77=7, 5\", \"47\u00f75=9, 2\"],\n  [\"89\u00f76=14, 5\", \"44\u00f73=14, 2\"],\n  [\"74\u00f75=14, 4\", \"91\u00f78=11, 3\"],\n  [\"87\u00f73=29, 0\", \"12\u00f74=3, 0\"],\n  [\"91\u00f76=15, 1\", \"22\u00f78=2, 6\"],\n  [\"61\u00f76=10, 1\", \"49\u00f76=8, 1\"],\n  [\"94\u00f75=18, 4\", \"75\u00f74=18, 3\"],\n  [\"41\u00f76=6, 5\", \"48\u00f72=24, 0\"],\n  [\"10\u00f77=1, 3\", \"94\u00f77=13, 3\"],\n  [\"78\u00f73=26, 0\", \"89\u00f72=44, 1\"],\n  [\"27\u00f75=5, 2\", \"83\u00f78=10, 3\"],\n  [\"13\u00f75=2, 3\", \"47\u00f75=9, 2\"],\n  [\"79\u00f78=9, 7\", \"10\u00f78=1, 2\"],\n  [\"54\u00f73=18, 0\", \"46\u00f75=9, 1\"],\n  [\"15\u00f78=1, 7\", \"72\u00f79=8, 0\"],\n  [\"70\u00f75=14, 0\", \"59\u00f73=19, 2\"],\n  [\"41\u00f78=5, 1\", \"81\u00f79=9, 0\"],\n  [\"73\u00f76=12, 1\", \"97\u00f75=19, 2\"],\n  [\"41\u00f77=5, 6\", \"26\u00f75=5, 1\"],\n  [\"89\u00f74=22, 1\", \"49\u00f73=16, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  // Each `oldText` is a full \"XX\u00f7Y=Q, R\" equation string unique to one\n  // table cell, so a plain case-sensitive substring search is enough\n  // (no need for matchWholeWord - the \"\u00f7\"/\",\"/\"=\" characters aren't\n  // word characters anyway).\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  // Each source string is unique in the document, so there should be\n  // exactly one hit - but replace every hit defensively in case a cell\n  // repeats the same value.\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit-divided-by-one-digit answer table: each table\n# cell holds a single \"XX\u00f7Y=Q, R\" style answer string that is being\n# replaced with a new equation/answer string (same cell, same\n# formatting - only the text content changes).\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"74\u00f78=9, 2\", \"37\u00f76=6, 1\"),\n    @(\"94\u00f78=11, 6\", \"75\u00f73=25, 0\"),\n    @(\"65\u00f72=32, 1\", \"84\u00f72=42, 0\"),\n    @(\"94\u00f76=15, 4\", \"36\u00f72=18, 0\"),\n    @(\"25\u00f73=8, 1\", \"86\u00f74=21, 2\"),\n    @(\"54\u00f77=7, 5\", \"47\u00f75=9, 2\"),\n    @(\"89\u00f76=14, 5\", \"44\u00f73=14, 2\"),\n    @(\"74\u00f75=14, 4\", \"91\u00f78=11, 3\"),\n    @(\"87\u00f73=29, 0\", \"12\u00f74=3, 0\"),\n    @(\"91\u00f76=15, 1\", \"22\u00f78=2, 6\"),\n    @(\"61\u00f76=10, 1\", \"49\u00f76=8, 1\"),\n    @(\"94\u00f75=18, 4\", \"75\u00f74=18, 3\"),\n    @(\"41\u00f76=6, 5\", \"48\u00f72=24, 0\"),\n    @(\"10\u00f77=1, 3\", \"94\u00f77=13, 3\"),\n    @(\"78\u00f73=26, 0\", \"89\u00f72=44, 1\"),\n    @(\"27\u00f75=5, 2\", \"83\u00f78=10, 3\"),\n    @(\"13\u00f75=2, 3\", \"47\u00f75=9, 2\"),\n    @(\"79\u00f78=9, 7\", \"10\u00f78=1, 2\"),\n    @(\"54\u00f73=18, 0\", \"46\u00f75=9, 1\"),\n    @(\"15\u00f78=1, 7\", \"72\u00f79=8, 0\"),\n    @(\"70\u00f75=14, 0\", \"59\u00f73=19, 2\"),\n    @(\"41\u00f78=5, 1\", \"81\u00f79=9, 0\"),\n    @(\"73\u00f76=12, 1\", \"97\u00f75=19, 2\"),\n    @(\"41\u00f77=5, 6\", \"26\u00f75=5, 1\"),\n    @(\"89\u00f74=22, 1\", \"49\u00f73=16, 1\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, \"wdFindContinue\", $false, $newText, \"wdReplaceAll\") | Out-Null\n}\n"}
